$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132, pushing existing rows 132-258 down to 133-259
$ws.Rows("132").Insert()

# Populate the newly inserted row 132 with the new daily price record
$ws.Range("A132").Value = 3
$ws.Range("B132").Value = "Femacal de La Calera"
$ws.Range("C132").Value = "Coquimbo"
$ws.Range("D132").Value = 45033
$ws.Range("E132").Value = 5
$ws.Range("F132").Value = 100112052
$ws.Range("G132").Value = "Albahaca"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 60
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = 4500
$ws.Range("N132").Value = "$/docena de matas"
$ws.Range("O132").Value = "Provincia de Quillota"
$ws.Range("P132").Value = 750
$ws.Range("Q132").Value = 6
$ws.Range("R132").Value = "Hortaliza"
